# Apply corrections to official place names, per commit message:
# "corrected most names to the official names from website"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# G3: Ramanagar -> Ramanagara
$ws.Range("G3").Value = "Ramanagara"

# G4: Ramanagar -> Ramanagara
$ws.Range("G4").Value = "Ramanagara"

# G12: Ramangara -> Ramanagara
$ws.Range("G12").Value = "Ramanagara"

# F17: remove the empty inline-string cell entirely (clear contents)
$ws.Range("F17").ClearContents()

# G33: Basavakalyan -> Bidar
$ws.Range("G33").Value = "Bidar"

# F34: remove the empty inline-string cell entirely (clear contents)
$ws.Range("F34").ClearContents()
